# JS_React_Syllabus.xlsx edit:
#  1. Insert a new topic row ("Version Control: Git and GitHub") right after the
#     "ClassWork: API Fetch" row, inside the "Day 4 / Introduction to React" block.
#  2. Make the "Day 3 / API and JavaScript Concepts" block's Day/Module Name
#     columns bold (matching the other bold day blocks).
#  3. Update the active selection to land on the newly inserted row's topic cell.

$wb = $excel.ActiveWorkbook
# "Sheet2" is the syllabus table and is the tab-selected / active sheet.
$ws = $wb.ActiveSheet

# --- 1. Insert the new row --------------------------------------------------
# Row 16 currently holds "Overview of React and its Architecture" (first topic
# of the Day 4 block after "ClassWork: API Fetch" on row 15). Inserting here
# shifts it (and everything below) down by one row and auto-extends the
# A15:A19 / B15:B19 merged "Day 4" / "Introduction to React" cells to
# A15:A20 / B15:B20.
$ws.Rows.Item(16).Insert() | Out-Null

# New row 16 inherits formatting from the row above; just fill in the topic.
$ws.Range("C16").Value = "Version Control: Git and GitHub"

# --- 2. Bold the Day 3 block (rows 10-14) -----------------------------------
$ws.Range("A10:B14").Font.Bold = $true

# --- 3. Update selection/scroll position ------------------------------------
$ws.Range("C17").Select() | Out-Null
$excel.ActiveWindow.ScrollRow = 7
